# The commit swaps the two theme parts shipped in the deck:
#   ppt/theme/theme1.xml  (the deck's "Integral" / Red Violet colour theme,
#                           used by the slide master) <->
#   ppt/theme/theme2.xml  (the "Office Theme" colours that used to only be
#                           referenced by the notes master)
# After the edit theme1.xml carries the standard "Office" palette.
#
# The font scheme and format scheme (fills/lines/effects) are byte-identical
# between the two theme parts, so the only externally visible difference is
# the colour scheme (and its/its parent's display name, which PowerPoint
# does not expose as a writable property through the object model). We
# therefore re-point every slot of the active theme's colour scheme -- the
# one seated on the slide master, which is what every slide in the deck
# actually renders with -- to the target "Office" RGB values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme

# PowerPoint COM RGB longs are 0x00BBGGRR, i.e. the reverse byte order of
# the usual RRGGBB hex notation used in OOXML's <a:srgbClr val=".."/>.
function Hex-ToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme (this is what ppt/theme/theme2.xml already
# holds before the edit), in ThemeColorScheme.Colors(1..12) order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Hex-ToRgb $officeColors[$i - 1]
}

# Best-effort: PowerPoint also labels the applied theme / colour scheme
# "Office Theme" / "Office" respectively once this palette is in place.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
